$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("H:H").Cut()
$ws.Columns("A:A").Insert()
$ws.Columns("L:L").Delete()
